$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value = 1000002
$ws.Range("C5").Value = "刀光-蓝色-旋转"
$ws.Range("D5").Value = "Effect_Sword_Slash_2"

$ws.Range("B6").Value = 1000003
$ws.Range("C6").Value = "刀光-浅蓝色-力量竖劈"
$ws.Range("D6").Value = "Effect_Sword_Slash_3"

$ws.Range("D7").Select()
